$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("E2")

# Populate the new analysis columns (éxito, absentismo_justificado,
# absentismo_injustificado, ratio) for every group row, and append the new
# "3DAM" (CFGS) row 38 that was added to the E2 sheet.
$ws.Range("C2").Value = 30.8
$ws.Range("D2").Value = 1.42
$ws.Range("E2").Value = 4.07
$ws.Range("H2").Value = 26
$ws.Range("C3").Value = 34.6
$ws.Range("D3").Value = 1.06
$ws.Range("E3").Value = 4.77
$ws.Range("H3").Value = 26
$ws.Range("C4").Value = 35.7
$ws.Range("D4").Value = 1.55
$ws.Range("E4").Value = 3.38
$ws.Range("H4").Value = 28
$ws.Range("C5").Value = 48.1
$ws.Range("D5").Value = 1.69
$ws.Range("E5").Value = 5.23
$ws.Range("H5").Value = 27
$ws.Range("C6").Value = 56
$ws.Range("D6").Value = 1.39
$ws.Range("E6").Value = 6.52
$ws.Range("H6").Value = 25
$ws.Range("C7").Value = 60
$ws.Range("D7").Value = 1.66
$ws.Range("E7").Value = 4.77
$ws.Range("H7").Value = 15
$ws.Range("C8").Value = 16.7
$ws.Range("D8").Value = 1.02
$ws.Range("E8").Value = 7.36
$ws.Range("H8").Value = 6
$ws.Range("C9").Value = 53.8
$ws.Range("D9").Value = 2.04
$ws.Range("E9").Value = 3.06
$ws.Range("H9").Value = 26
$ws.Range("C10").Value = 30.8
$ws.Range("D10").Value = 2.36
$ws.Range("E10").Value = 4.74
$ws.Range("H10").Value = 26
$ws.Range("C11").Value = 39.3
$ws.Range("D11").Value = 2.82
$ws.Range("E11").Value = 2.75
$ws.Range("H11").Value = 28
$ws.Range("C12").Value = 25
$ws.Range("D12").Value = 3.65
$ws.Range("E12").Value = 4.85
$ws.Range("H12").Value = 12
$ws.Range("C13").Value = 45.8
$ws.Range("D13").Value = 4.13
$ws.Range("E13").Value = 2.55
$ws.Range("H13").Value = 24
$ws.Range("C14").Value = 33.3
$ws.Range("D14").Value = 3.01
$ws.Range("E14").Value = 3.46
$ws.Range("H14").Value = 27
$ws.Range("C15").Value = 24
$ws.Range("D15").Value = 3.69
$ws.Range("E15").Value = 6.99
$ws.Range("H15").Value = 25
$ws.Range("C16").Value = 46.7
$ws.Range("D16").Value = 5.53
$ws.Range("E16").Value = 20.56
$ws.Range("H16").Value = 15
$ws.Range("C17").Value = 14.3
$ws.Range("D17").Value = 8.67
$ws.Range("E17").Value = 20.87
$ws.Range("H17").Value = 7
$ws.Range("C18").Value = 29.4
$ws.Range("D18").Value = 3.02
$ws.Range("E18").Value = 2.74
$ws.Range("H18").Value = 17
$ws.Range("C19").Value = 25
$ws.Range("D19").Value = 2.52
$ws.Range("E19").Value = 6.56
$ws.Range("H19").Value = 24
$ws.Range("C20").Value = 52.6
$ws.Range("D20").Value = 2.17
$ws.Range("E20").Value = 3.7
$ws.Range("H20").Value = 19
$ws.Range("C21").Value = 29.2
$ws.Range("D21").Value = 2.49
$ws.Range("E21").Value = 4.65
$ws.Range("H21").Value = 24
$ws.Range("C22").Value = 40
$ws.Range("D22").Value = 3.72
$ws.Range("E22").Value = 19.05
$ws.Range("H22").Value = 15
$ws.Range("C23").Value = 31.2
$ws.Range("D23").Value = 2.01
$ws.Range("E23").Value = 15.08
$ws.Range("H23").Value = 16
$ws.Range("C24").Value = 57.9
$ws.Range("D24").Value = 3.49
$ws.Range("E24").Value = 17.23
$ws.Range("H24").Value = 19
$ws.Range("C25").Value = 100
$ws.Range("D25").Value = 3.49
$ws.Range("E25").Value = 14
$ws.Range("H25").Value = 8
$ws.Range("C26").Value = 61.5
$ws.Range("D26").Value = 1.7
$ws.Range("E26").Value = 11.06
$ws.Range("H26").Value = 13
$ws.Range("C27").Value = 81.8
$ws.Range("D27").Value = 1.59
$ws.Range("E27").Value = 8.69
$ws.Range("H27").Value = 11
$ws.Range("C28").Value = 66.7
$ws.Range("D28").Value = 1.81
$ws.Range("E28").Value = 6.42
$ws.Range("H28").Value = 12
$ws.Range("C29").Value = 74.09999999999999
$ws.Range("D29").Value = 1.97
$ws.Range("E29").Value = 11.54
$ws.Range("H29").Value = 27
$ws.Range("C30").Value = 67.90000000000001
$ws.Range("D30").Value = 4.27
$ws.Range("E30").Value = 10.48
$ws.Range("H30").Value = 28
$ws.Range("C31").Value = 57.7
$ws.Range("D31").Value = 1.48
$ws.Range("E31").Value = 11.84
$ws.Range("H31").Value = 26
$ws.Range("C32").Value = 28
$ws.Range("D32").Value = 3.14
$ws.Range("E32").Value = 6.93
$ws.Range("H32").Value = 25
$ws.Range("C33").Value = 93.3
$ws.Range("D33").Value = 1.12
$ws.Range("E33").Value = 4.56
$ws.Range("H33").Value = 15
$ws.Range("C34").Value = 100
$ws.Range("D34").Value = 1.51
$ws.Range("E34").Value = 5.51
$ws.Range("H34").Value = 23
$ws.Range("C35").Value = 90.90000000000001
$ws.Range("D35").Value = 1.95
$ws.Range("E35").Value = 7.58
$ws.Range("H35").Value = 22
$ws.Range("C36").Value = 94.7
$ws.Range("D36").Value = 1.46
$ws.Range("E36").Value = 5.62
$ws.Range("H36").Value = 19
$ws.Range("C37").Value = 58.3
$ws.Range("D37").Value = 1.68
$ws.Range("E37").Value = 16.29
$ws.Range("H37").Value = 12
$ws.Range("A38").Value = "3DAM"
$ws.Range("B38").Value = "CFGS"
$ws.Range("C38").Value = 100
$ws.Range("D38").Value = 0
$ws.Range("E38").Value = 1.01
$ws.Range("H38").Value = 5

# Set up the view state for E3 first (it loses the active/tabSelected tab).
$ws3 = $wb.Worksheets.Item("E3")
$ws3.Activate()
$ws3.Range("G18").Select()

# Then activate E2 last, making it the active sheet (workbookView activeTab=1)
# with its own new selection reflecting the freshly added row 38.
$ws.Activate()
$ws.Range("C38").Select()
